# Scheduled runner update: refresh market-price derived columns (H-N)
# across the leve-profit tables in each job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 366.66666
$ws.Range("I12").Value = 366.66666
$ws.Range("K12").Value = 366.66666
$ws.Range("M12").Value = -196.66666

# Row 17
$ws.Range("H17").Value = 424.57693
$ws.Range("J17").Value = 424.57693
$ws.Range("L17").Value = 1273.73079
$ws.Range("N17").Value = -1609.73079

# Row 33
$ws.Range("H33").Value = 157.45454
$ws.Range("I33").Value = 81.55556
$ws.Range("K33").Value = 81.55556
$ws.Range("M33").Value = 147.44444

# Row 43
$ws.Range("H43").Value = 27797778
$ws.Range("I43").Value = 40001
$ws.Range("J43").Value = 55555556
$ws.Range("K43").Value = 40001
$ws.Range("L43").Value = 55555556
$ws.Range("M43").Value = -39932
$ws.Range("N43").Value = -55555694

# Row 98
$ws.Range("H98").Value = 2414.8572
$ws.Range("I98").Value = 2006.2222
$ws.Range("J98").Value = 4866.6665
$ws.Range("K98").Value = 2006.2222
$ws.Range("L98").Value = 4866.6665
$ws.Range("M98").Value = -508.2221999999999
$ws.Range("N98").Value = -7862.6665

# Row 116
$ws.Range("H116").Value = 2836.1292
$ws.Range("I116").Value = 2490.2632
$ws.Range("J116").Value = 3383.75
$ws.Range("K116").Value = 2490.2632
$ws.Range("L116").Value = 3383.75
$ws.Range("M116").Value = 951.7368000000001
$ws.Range("N116").Value = -10267.75

# Row 122
$ws.Range("H122").Value = 2414.8572
$ws.Range("I122").Value = 2006.2222
$ws.Range("J122").Value = 4866.6665
$ws.Range("K122").Value = 6018.6666
$ws.Range("L122").Value = 14599.9995
$ws.Range("M122").Value = -3568.6666
$ws.Range("N122").Value = -19499.9995

# Row 137
$ws.Range("H137").Value = 1098.2941
$ws.Range("I137").Value = 1008.069
$ws.Range("J137").Value = 1621.6
$ws.Range("K137").Value = 3024.207
$ws.Range("L137").Value = 4864.799999999999
$ws.Range("M137").Value = -474.2069999999999
$ws.Range("N137").Value = -9964.799999999999

$ws = $wb.Worksheets.Item("ARM")
# Row 31
$ws.Range("H31").Value = 19166
$ws.Range("I31").Value = 4165.3335
$ws.Range("J31").Value = 34166.668
$ws.Range("K31").Value = 4165.3335
$ws.Range("L31").Value = 34166.668
$ws.Range("M31").Value = -3871.3335
$ws.Range("N31").Value = -34754.668

# Row 32
$ws.Range("H32").Value = 3998.745
$ws.Range("I32").Value = 3766.9092
$ws.Range("K32").Value = 3766.9092
$ws.Range("M32").Value = -3479.9092

# Row 45
$ws.Range("H45").Value = 1450.9333
$ws.Range("I45").Value = 1519.6923
$ws.Range("J45").Value = 1004
$ws.Range("K45").Value = 1519.6923
$ws.Range("L45").Value = 1004
$ws.Range("M45").Value = -1142.6923
$ws.Range("N45").Value = -1758

# Row 107
$ws.Range("H107").Value = 16714
$ws.Range("J107").Value = 16714
$ws.Range("L107").Value = 16714
$ws.Range("N107").Value = -24394

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 90910216
$ws.Range("I16").Value = 100001090
$ws.Range("K16").Value = 100001090
$ws.Range("M16").Value = -100000803

# Row 31
$ws.Range("H31").Value = 2099.7827
$ws.Range("I31").Value = 1117.7273
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1117.7273
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -822.7273
$ws.Range("N31").Value = -3590

# Row 34
$ws.Range("H34").Value = 2099.7827
$ws.Range("I34").Value = 1117.7273
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1117.7273
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -915.7273
$ws.Range("N34").Value = -3404

# Row 99
$ws.Range("H99").Value = 1779.9524
$ws.Range("I99").Value = 1740.7142
$ws.Range("J99").Value = 1858.4286
$ws.Range("K99").Value = 1740.7142
$ws.Range("L99").Value = 1858.4286
$ws.Range("M99").Value = -242.7141999999999
$ws.Range("N99").Value = -4854.4286

# Row 107
$ws.Range("H107").Value = 526.4194
$ws.Range("I107").Value = 420.57144
$ws.Range("K107").Value = 420.57144
$ws.Range("M107").Value = 1499.42856

# Row 113
$ws.Range("H113").Value = 90910216
$ws.Range("I113").Value = 100001090
$ws.Range("K113").Value = 100001090
$ws.Range("M113").Value = -99998920

# Row 126
$ws.Range("H126").Value = 1779.9524
$ws.Range("I126").Value = 1740.7142
$ws.Range("J126").Value = 1858.4286
$ws.Range("K126").Value = 5222.142599999999
$ws.Range("L126").Value = 5575.2858
$ws.Range("M126").Value = -2752.142599999999
$ws.Range("N126").Value = -10515.2858

# Row 132
$ws.Range("H132").Value = 9813.571
$ws.Range("I132").Value = 11989.2
$ws.Range("J132").Value = 4374.5
$ws.Range("K132").Value = 35967.60000000001
$ws.Range("L132").Value = 13123.5
$ws.Range("M132").Value = -33437.60000000001
$ws.Range("N132").Value = -18183.5

# Row 134
$ws.Range("H134").Value = 15874693
$ws.Range("I134").Value = 22223918
$ws.Range("J134").Value = 1627.8334
$ws.Range("K134").Value = 66671754
$ws.Range("L134").Value = 4883.5002
$ws.Range("M134").Value = -66669219
$ws.Range("N134").Value = -9953.5002

$ws = $wb.Worksheets.Item("CUL")
# Row 118
$ws.Range("H118").Value = 489.66666
$ws.Range("I118").Value = 489.66666
$ws.Range("K118").Value = 1468.99998
$ws.Range("M118").Value = -225.9999800000001

# Row 131
$ws.Range("H131").Value = 14707129
$ws.Range("I131").Value = 142857440
$ws.Range("J131").Value = 1355.459
$ws.Range("K131").Value = 428572320
$ws.Range("L131").Value = 4066.377
$ws.Range("M131").Value = -428567280
$ws.Range("N131").Value = -14146.377

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 500
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 500
$ws.Range("N2").Value = -726
$ws.Range("M2").ClearContents()

# Row 122
$ws.Range("H122").Value = 108385.64
$ws.Range("I122").Value = 1324.9166
$ws.Range("K122").Value = 3974.7498
$ws.Range("M122").Value = -1524.7498

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2642.0833
$ws.Range("I40").Value = 2533.3333
$ws.Range("K40").Value = 2533.3333
$ws.Range("M40").Value = -2397.3333

# Row 46
$ws.Range("H46").Value = 3898.6
$ws.Range("I46").Value = 1828.6666
$ws.Range("J46").Value = 4785.7144
$ws.Range("K46").Value = 1828.6666
$ws.Range("L46").Value = 4785.7144
$ws.Range("M46").Value = -1640.6666
$ws.Range("N46").Value = -5161.7144

# Row 122
$ws.Range("H122").Value = 22730398
$ws.Range("I122").Value = 31253174
$ws.Range("K122").Value = 93759522
$ws.Range("M122").Value = -93757072

# Row 136
$ws.Range("H136").Value = 13076.333
$ws.Range("I136").Value = 21637.6
$ws.Range("J136").Value = 2374.75
$ws.Range("K136").Value = 64912.8
$ws.Range("L136").Value = 7124.25
$ws.Range("M136").Value = -62362.8
$ws.Range("N136").Value = -12224.25

# Row 139
$ws.Range("H139").Value = 50665
$ws.Range("J139").Value = 50665
$ws.Range("L139").Value = 50665
$ws.Range("N139").Value = -60945

# Row 140
$ws.Range("H140").Value = 46891.125
$ws.Range("J140").Value = 46891.125
$ws.Range("L140").Value = 46891.125
$ws.Range("N140").Value = -57251.125

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 9287567
$ws.Range("I122").Value = 9287567
$ws.Range("K122").Value = 27862701
$ws.Range("M122").Value = -27860251

# Row 136
$ws.Range("H136").Value = 1162.6364
$ws.Range("I136").Value = 976.6667
$ws.Range("K136").Value = 2930.0001
$ws.Range("M136").Value = -380.0001000000002
